$wb = $excel.ActiveWorkbook

# The data that changed lives on the worksheet named "100"
$ws = $wb.Worksheets.Item("100")

$ws.Range("B3").Value = 92.80000305175781
$ws.Range("B4").Value = 87.09999847412109
$ws.Range("B5").Value = 85.69999694824219
$ws.Range("B6").Value = 58.56666564941406
$ws.Range("B7").Value = 57.59999847412109
$ws.Range("B8").Value = 51.79999923706055
$ws.Range("B9").Value = 50.83333206176758
$ws.Range("B10").Value = 26.41666603088379
$ws.Range("B11").Value = 25.39999961853027
